# Applies the "Updated cryptos list" data refresh described in the commit diff.
# Each Range.Value assignment mirrors one <t> text change from the OOXML diff.
# Numeric-looking price strings are prefixed with a literal apostrophe so Excel
# stores them as text (preserving trailing zeros / multi-dot formatting) instead
# of silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.061.38"
$ws.Range("E2").Value = "  +0.88%  "

$ws.Range("D3").Value = "2.211.14"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").Value = "'231.05"
$ws.Range("E5").Value = "  +1.17%  "

$ws.Range("D6").Value = "'0.617"
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("D7").Value = "'60.90"
$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("D9").Value = "'0.402"
$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").Value = "'0.0898"
$ws.Range("E10").Value = "  +2.75%  "

$ws.Range("E11").Value = "  -0.50%  "

$ws.Range("D12").Value = "2.534.38"
$ws.Range("E12").Value = "  -0.93%  "

$ws.Range("D13").Value = "'15.43"
$ws.Range("E13").Value = "  -1.41%  "

$ws.Range("D14").Value = "'21.97"
$ws.Range("E14").Value = "  +2.62%  "

$ws.Range("D15").Value = "'0.796"
$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("E16").Value = "  +0.47%  "

$ws.Range("D17").Value = "2.213.97"
$ws.Range("E17").Value = "  -0.64%  "

$ws.Range("D18").Value = "41.918.67"
$ws.Range("E18").Value = "  +0.80%  "

$ws.Range("E19").Value = "  +5.12%  "

$ws.Range("D20").Value = "'6.20"
$ws.Range("E20").Value = "  +2.95%  "

$ws.Range("D21").Value = "'71.86"
$ws.Range("E21").Value = "  -0.90%  "

$ws.Range("D22").Value = "'242.73"
$ws.Range("E22").Value = "  -1.61%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").Value = "'2.43"
$ws.Range("E24").Value = "  +3.22%  "

$ws.Range("D25").Value = "'2.29"
$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("D26").Value = "'9.54"
$ws.Range("E26").Value = "  -0.57%  "

$ws.Range("D27").Value = "'168.86"
$ws.Range("E27").Value = "  +0.69%  "

$ws.Range("E28").Value = "  +0.70%  "

$ws.Range("D29").Value = "'20.21"
$ws.Range("E29").Value = "  +1.75%  "

$ws.Range("D30").Value = "'1.45"
$ws.Range("E30").Value = "  +2.08%  "

$ws.Range("D31").Value = "'2.66"
$ws.Range("E31").Value = "  +2.82%  "

$ws.Range("E32").Value = "  -1.29%  "

$ws.Range("D33").Value = "'4.95"
$ws.Range("E33").Value = "  -1.79%  "

$ws.Range("D34").Value = "'4.59"
$ws.Range("E34").Value = "  -1.13%  "

$ws.Range("D35").Value = "'0.0646"
$ws.Range("E35").Value = "  +3.95%  "

$ws.Range("D36").Value = "'6.29"
$ws.Range("E36").Value = "  -4.60%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'3.53"
$ws.Range("E37").Value = "  -4.42%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.33"
$ws.Range("E38").Value = "  -1.37%  "

$ws.Range("D39").Value = "'0.0249"
$ws.Range("E39").Value = "  +6.41%  "

$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("B41").Value = "TerraClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D41").Value = "'0.000226"
$ws.Range("E41").Value = "  -5.15%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'8.54"
$ws.Range("E42").Value = "  -2.88%  "

$ws.Range("D43").Value = "'0.0952"
$ws.Range("E43").Value = "  -2.16%  "

$ws.Range("D44").Value = "'1.20"
$ws.Range("E44").Value = "  +1.44%  "

$ws.Range("D45").Value = "'96.76"
$ws.Range("E45").Value = "  -2.26%  "

$ws.Range("D46").Value = "1.456.90"
$ws.Range("E46").Value = "  -0.59%  "

$ws.Range("D47").Value = "'4.26"
$ws.Range("E47").Value = "  -12.45%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'16.10"
$ws.Range("E48").Value = "  -1.11%  "

$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "'2.72"
$ws.Range("E49").Value = "  -2.24%  "

$ws.Range("E50").Value = "  -1.10%  "

$ws.Range("D51").Value = "'2.20"
$ws.Range("E51").Value = "  +1.90%  "
